# Updates cryptos list price/volume figures (and restores the FTXToken/ARBITRUM
# row ordering) to match the latest scrape, per the "Updated cryptos list" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = '37.785.34'
$ws.Range("E2").Value = '  -0.09%  '

# Row 3 - Ethereum
$ws.Range("D3").Value = '2.075.95'
$ws.Range("E3").Value = '  -0.51%  '

# Row 4 - TetherUSD
$ws.Range("E4").Value = '  +0.11%  '

# Row 5 - BNB
$ws.Range("E5").Value = '  -0.70%  '

# Row 6 - XRP
$ws.Range("D6").Value = '''0.624'
$ws.Range("E6").Value = '  -0.18%  '

# Row 7 - Solana
$ws.Range("D7").Value = '''58.43'
$ws.Range("E7").Value = '  -2.12%  '

# Row 9 - Cardano
$ws.Range("E9").Value = '  +0.22%  '

# Row 10 - Dogecoin
$ws.Range("D10").Value = '''0.0784'
$ws.Range("E10").Value = '  -1.06%  '

# Row 11 - TRON
$ws.Range("E11").Value = '  +3.82%  '

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = '2.381.89'
$ws.Range("E12").Value = '  -0.40%  '

# Row 13 - Chainlink
$ws.Range("D13").Value = '''14.76'
$ws.Range("E13").Value = '  +0.41%  '

# Row 14 - Avalanche
$ws.Range("D14").Value = '''21.12'
$ws.Range("E14").Value = '  -1.77%  '

# Row 15 - Polygon
$ws.Range("E15").Value = '  +0.63%  '

# Row 16 - Polkadot
$ws.Range("D16").Value = '''5.33'
$ws.Range("E16").Value = '  -0.05%  '

# Row 17 - WrappedEther
$ws.Range("D17").Value = '2.083.42'
$ws.Range("E17").Value = '  -0.19%  '

# Row 18 - WrappedBTC
$ws.Range("D18").Value = '37.664.11'
$ws.Range("E18").Value = '  -0.21%  '

# Row 19 - Uniswap
$ws.Range("D19").Value = '''6.13'
$ws.Range("E19").Value = '  -1.29%  '

# Row 20 - Litecoin
$ws.Range("D20").Value = '''71.47'
$ws.Range("E20").Value = '  -0.31%  '

# Row 21 - ShibaInu
$ws.Range("D21").Value = '0.0₃0841'
$ws.Range("E21").Value = '  +1.26%  '

# Row 22 - BitcoinCash
$ws.Range("D22").Value = '''228.99'
$ws.Range("E22").Value = '  +0.02%  '

# Row 23 - Dai
$ws.Range("D23").Value = '''0.999'
$ws.Range("E23").Value = '  -0.08%  '

# Row 24 - Toncoin
$ws.Range("D24").Value = '''2.40'
$ws.Range("E24").Value = '  -0.66%  '

# Row 25 - PancakeSwap
$ws.Range("E25").Value = '  -0.61%  '

# Row 26 - Cosmos
$ws.Range("E26").Value = '  +7.26%  '

# Row 27 - Monero
$ws.Range("D27").Value = '''171.74'
$ws.Range("E27").Value = '  +0.75%  '

# Row 28 - Kaspa
$ws.Range("D28").Value = '''0.139'
$ws.Range("E28").Value = '  -0.22%  '

# Row 29 - ImmutableX
$ws.Range("E29").Value = '  -2.45%  '

# Row 30 - EthereumClassic
$ws.Range("D30").Value = '''19.42'
$ws.Range("E30").Value = '  -0.69%  '

# Row 31 - Stellar
$ws.Range("E31").Value = '  +1.29%  '

# Row 32 - Filecoin
$ws.Range("D32").Value = '''4.73'
$ws.Range("E32").Value = '  +0.11%  '

# Row 33 - Hedera
$ws.Range("E33").Value = '  +0.24%  '

# Row 34 - InternetComputer(DFINITY)
$ws.Range("D34").Value = '''4.71'
$ws.Range("E34").Value = '  +0.22%  '

# Row 35 - LidoDAOToken
$ws.Range("E35").Value = '  -2.88%  '

# Row 36 - WEMIXToken
$ws.Range("D36").Value = '''1.82'
$ws.Range("E36").Value = '  -0.18%  '

# Row 37 - RenderToken
$ws.Range("D37").Value = '''3.40'
$ws.Range("E37").Value = '  -3.75%  '

# Row 38 - BinanceUSD
$ws.Range("D38").Value = '''0.999'
$ws.Range("E38").Value = '  -0.01%  '

# Row 39 - THORChain
$ws.Range("D39").Value = '''5.44'
$ws.Range("E39").Value = '  -0.19%  '

# Row 40 - VeChain
$ws.Range("E40").Value = '  +8.46%  '

# Row 41 - Aave
$ws.Range("D41").Value = '''100.20'
$ws.Range("E41").Value = '  -0.17%  '

# Row 42 - Cronos
$ws.Range("E42").Value = '  -1.32%  '

# Row 43 - InjectiveProtocol
$ws.Range("D43").Value = '''17.20'
$ws.Range("E43").Value = '  +6.18%  '

# Row 44 - HuobiToken
$ws.Range("E44").Value = '  -0.99%  '

# Row 45 - Maker
$ws.Range("D45").Value = '1.446.98'
$ws.Range("E45").Value = '  -1.28%  '

# Row 46 - TrustWalletToken
$ws.Range("D46").Value = '''1.15'
$ws.Range("E46").Value = '  -1.79%  '

# Row 47 - FTXToken
$ws.Range("B47").Value = 'ARBITRUM'
$ws.Range("C47").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D47").Value = '''1.06'
$ws.Range("E47").Value = '  -0.39%  '

# Row 48 - ARBITRUM
$ws.Range("B48").Value = 'FTXToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D48").Value = '''4.11'
$ws.Range("E48").Value = '  -2.86%  '

# Row 49 - FraxShare
$ws.Range("D49").Value = '''7.38'
$ws.Range("E49").Value = '  -1.28%  '

# Row 50 - MXToken
$ws.Range("E50").Value = '  -1.72%  '

# Row 51 - RocketPoolETH
$ws.Range("D51").Value = '2.267.52'
$ws.Range("E51").Value = '  -0.42%  '
